$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "As you know, the PGB approved the Naming & Design Rules (NDR)
#    v6.0 PSD01 and the NIEM Model Version v6.0 PS02 will be out for
#    vote soon."
#    ->
#    "As you know, the PGB is voting on the Naming & Design Rules
#    (NDR) v6.0 PSD01. There  are currently 9 votes to approve which
#    is a majority. The NIEM Model Version v6.0 PS02 will be out for
#    vote soon."
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "the PGB approved the Naming & Design Rules (NDR) v6.0 PSD01 and the NIEM Model Version v6.0 PS02 will",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the PGB is voting on the Naming & Design Rules (NDR) v6.0 PSD01. There  are currently 9 votes to approve which is a majority. The NIEM Model Version v6.0 PS02 will",
    2) | Out-Null

# ------------------------------------------------------------------
# 2) "At our meeting on 30 Jan, we expect to motion:"
#    -> "At our meeting on 30 Jan, we expect to motion to:"
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "we expect to motion:", $true, $false, $false, $false, $false, $true, 1, $false,
    "we expect to motion to:", 2) | Out-Null

# ------------------------------------------------------------------
# 3) First bullet:
#    "Advance NDR v6.0 PSD01 to a project specification (PS), and"
#    -> "Approve the 14 Nov PGB Meeting Draft Minutes"
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Advance NDR v6.0 PSD01 to a project specification (PS), and",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Approve the 14 Nov PGB Meeting Draft Minutes", 2) | Out-Null

# ------------------------------------------------------------------
# 4) Second bullet:
#    "Intent to advance NIEM Model Version v6.0 PS02 to and OASIS
#    Standard (OS) once PS02 is approved."
#    -> split into two bullets:
#       "Update 2025 meeting dates to address scheduling conflicts"
#       "Advance the NIEM Model Version v6.0 PS02 to and OASIS
#        Standard (OS) once PS02 is approved."
#    and a brand new bullet inserted between them:
#       "Advance the approved draft of NDR v6.0 PSD01 to a Project
#        Specification (PS), and"
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Intent to advance NIEM Model Version v6.0 PS02 to and OASIS Standard (OS) once PS02 is approved.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Update 2025 meeting dates to address scheduling conflicts", 2) | Out-Null

# find the paragraph that now reads "Update 2025 meeting dates..."
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $par = $d.Paragraphs($i)
    if ($par.Range.Text -like "Update 2025 meeting dates to address scheduling conflicts*") {
        $target = $par
        break
    }
}

$target.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs($target.Range.ListFormat.ListValue -as [int])
# recompute index robustly: find the paragraph after $target
$idx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Start -eq $target.Range.Start) {
        $idx = $i
        break
    }
}
$newPara1 = $d.Paragraphs($idx + 1)
$newPara1.Range.Text = "Advance the approved draft of NDR v6.0 PSD01 to a Project Specification (PS), and"

$newPara1.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs($idx + 2)
$newPara2.Range.Text = "Advance the NIEM Model Version v6.0 PS02 to and OASIS Standard (OS) once PS02 is approved."

# ------------------------------------------------------------------
# 5) Add a blank paragraph plus the hackathon info-brief paragraph
#    right after the (now last) NIEM bullet, and before the blank
#    line that precedes "Microsoft Teams NIEMOpen PGB Meeting..."
# ------------------------------------------------------------------
$newPara2.Range.InsertParagraphAfter()
$blankPara = $d.Paragraphs($idx + 3)
$blankPara.Range.Text = ""

$blankPara.Range.InsertParagraphAfter()
$infoPara = $d.Paragraphs($idx + 4)
$infoPara.Range.Text = "The PGB will receive an info brief on hackathons as well as updates from NTAC, NBAC and NMO Technical Steering Committees."

